$d = $word.ActiveDocument

# --- Paragraph 1 (title) : drop w:hint="eastAsia" from the paragraph-mark rPr ---
$p1 = $d.Paragraphs.Item(1)
$p1Range = $p1.Range
$p1Xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="21D30246" w14:textId="5D76062D" w:rsidR="001D3745" w:rsidRPr="00AA10CA" w:rsidRDefault="00AA10CA"><w:pPr><w:spacing w:after="0" w:line="360" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="16"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" w:hint="eastAsia"/><w:b/><w:color w:val="000000"/><w:sz w:val="32"/></w:rPr><w:t>Robust</w:t></w:r><w:r w:rsidR="00000000"><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:color w:val="000000"/><w:sz w:val="32"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" w:hint="eastAsia"/><w:b/><w:color w:val="000000"/><w:sz w:val="32"/></w:rPr><w:t>estimations from distribution structures</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$p1Range.InsertXML($p1Xml)

# --- Paragraph 3 (authors line) : drop w:hint="eastAsia" from the paragraph-mark rPr ---
$p3 = $d.Paragraphs.Item(3)
$p3Range = $p3.Range
$p3Xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="3A134A3F" w14:textId="2B82CE34" w:rsidR="001D3745" w:rsidRPr="00AA10CA" w:rsidRDefault="00592EC5"><w:pPr><w:widowControl w:val="0"/><w:spacing w:after="0" w:line="360" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:i/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" w:hint="eastAsia"/><w:b/><w:i/><w:color w:val="000000"/><w:sz w:val="24"/><w:u w:val="single"/></w:rPr><w:t xml:space="preserve">Johon Li </w:t></w:r><w:r w:rsidR="00000000"><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:i/><w:color w:val="000000"/><w:sz w:val="24"/><w:u w:val="single"/></w:rPr><w:t>Tuobang</w:t></w:r><w:r w:rsidR="00000000"><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:i/><w:color w:val="000000"/><w:sz w:val="24"/><w:vertAlign w:val="superscript"/></w:rPr><w:t>1,2</w:t></w:r><w:r w:rsidR="00AA10CA"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" w:hint="eastAsia"/><w:i/><w:color w:val="000000"/><w:sz w:val="24"/><w:vertAlign w:val="superscript"/></w:rPr><w:t>,3</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$p3Range.InsertXML($p3Xml)

# --- Paragraph 9 (abstract body) : merge "sensitive" runs + split "discover..." run,
#     and drop w:hint="eastAsia" from the paragraph-mark rPr ---
$p9 = $d.Paragraphs.Item(9)
$p9Range = $p9.Range
$p9Xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="1475A954" w14:textId="014FF39C" w:rsidR="00592EC5" w:rsidRDefault="005029C3" w:rsidP="005029C3"><w:pPr><w:spacing w:after="200" w:line="360" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="005029C3"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/></w:rPr><w:t>Descriptive statistics for parametric</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" w:hint="eastAsia"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> or nonparametric</w:t></w:r><w:r w:rsidRPr="005029C3"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> models are </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" w:hint="eastAsia"/><w:sz w:val="24"/></w:rPr><w:t>generally</w:t></w:r><w:r w:rsidRPr="005029C3"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> sensitive to departures, gross errors, and/or random errors</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" w:hint="eastAsia"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">. Here, we explored semiparametric methods to </w:t></w:r><w:r w:rsidR="00741CB8"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" w:hint="eastAsia"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">classify distributions to </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" w:hint="eastAsia"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">reveal </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" w:hint="eastAsia"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">the underlying mechanisms </w:t></w:r><w:r w:rsidR="00741CB8"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" w:hint="eastAsia"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">of current robust estimators. </w:t></w:r><w:r w:rsidR="00741CB8" w:rsidRPr="00741CB8"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">Further deductions explain why the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00741CB8" w:rsidRPr="00741CB8"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/></w:rPr><w:t>Winsorized</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00741CB8" w:rsidRPr="00741CB8"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> mean typically has smaller biases compared to the trimmed mean</w:t></w:r><w:r w:rsidR="00741CB8"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" w:hint="eastAsia"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> and why the Hodges-Lehmann estimator </w:t></w:r><w:r w:rsidR="00741CB8"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/></w:rPr><w:t>and</w:t></w:r><w:r w:rsidR="00741CB8"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" w:hint="eastAsia"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> Bickel-Lehmann spread are the optimal nonparametric location and scale estimator. </w:t></w:r><w:r w:rsidR="000A6955"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" w:hint="eastAsia"/><w:sz w:val="24"/></w:rPr><w:t>From</w:t></w:r><w:r w:rsidR="00741CB8"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" w:hint="eastAsia"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> the</w:t></w:r><w:r w:rsidR="000A6955"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" w:hint="eastAsia"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> distribution structures, a series of new estimators were deduced. Some of them are </w:t></w:r><w:r w:rsidR="000A6955" w:rsidRPr="000A6955"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/></w:rPr><w:t>robust to both gross errors and departures from parametric assumptions, making them ideal for estimating the mean and central moments of common unimodal distributions.</w:t></w:r><w:r w:rsidR="005F4FF9"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" w:hint="eastAsia"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> This presentation </w:t></w:r><w:r w:rsidR="005F4FF9" w:rsidRPr="005F4FF9"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/></w:rPr><w:t>sheds light on the understanding of the common nature of probability distributions</w:t></w:r><w:r w:rsidR="005F4FF9"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" w:hint="eastAsia"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> and the measures of them.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$p9Range.InsertXML($p9Xml)

Write-Host "Edits applied"
